$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.751.32"
$ws.Range("E2").Value = "  -4.00%  "

$ws.Range("D3").Value = "3.351.52"
$ws.Range("E3").Value = "  -0.21%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'574.99"
$ws.Range("E5").Value = "  -2.74%  "

$ws.Range("D6").Value = "'182.30"
$ws.Range("E6").Value = "  -4.84%  "

$ws.Range("E7").Value = "  +2.93%  "

$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("D9").Value = "'0.129"
$ws.Range("E9").Value = "  -3.22%  "

$ws.Range("E10").Value = "  -1.02%  "

$ws.Range("D11").Value = "'0.405"
$ws.Range("E11").Value = "  -2.96%  "

$ws.Range("D12").Value = "3.935.74"
$ws.Range("E12").Value = "  -0.28%  "

$ws.Range("E13").Value = "  -0.75%  "

$ws.Range("D14").Value = "'26.94"
$ws.Range("E14").Value = "  -5.44%  "

$ws.Range("D15").Value = "66.877.95"
$ws.Range("E15").Value = "  -3.89%  "

$ws.Range("E16").Value = "  -1.95%  "

$ws.Range("D17").Value = "3.343.33"
$ws.Range("E17").Value = "  -1.00%  "

$ws.Range("D18").Value = "'435.90"
$ws.Range("E18").Value = "  -3.85%  "

$ws.Range("D19").Value = "'13.64"
$ws.Range("E19").Value = "  -0.56%  "

$ws.Range("D20").Value = "'5.69"
$ws.Range("E20").Value = "  -2.19%  "

$ws.Range("D21").Value = "'7.61"
$ws.Range("E21").Value = "  -3.08%  "

$ws.Range("D22").Value = "'73.62"
$ws.Range("E22").Value = "  -3.38%  "

$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.20%  "

$ws.Range("E24").Value = "  -0.56%  "

$ws.Range("D25").Value = "'0.0000117"
$ws.Range("E25").Value = "  -3.75%  "

$ws.Range("E26").Value = "  +1.28%  "

$ws.Range("D27").Value = "'9.08"
$ws.Range("E27").Value = "  -3.79%  "

$ws.Range("E28").Value = "  +0.04%  "

$ws.Range("D29").Value = "'1.96"
$ws.Range("E29").Value = "  -2.63%  "

$ws.Range("D30").Value = "'22.88"
$ws.Range("E30").Value = "  -1.74%  "

$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").Value = "'5.29"
$ws.Range("E31").Value = "  -4.57%  "

$ws.Range("B32").Value = "USDe"
$ws.Range("C32").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D32").Value = "'0.999"
$ws.Range("E32").Value = "  +0.06%  "

$ws.Range("D33").Value = "'6.80"
$ws.Range("E33").Value = "  -2.23%  "

$ws.Range("E34").Value = "  -4.87%  "

$ws.Range("D35").Value = "'160.64"
$ws.Range("E35").Value = "  -2.13%  "

$ws.Range("E36").Value = "  -5.03%  "

$ws.Range("D37").Value = "'28.03"
$ws.Range("E37").Value = "  +3.05%  "

$ws.Range("D38").Value = "'1.80"
$ws.Range("E38").Value = "  -7.38%  "

$ws.Range("D39").Value = "2.817.62"
$ws.Range("E39").Value = "  +3.69%  "

$ws.Range("D40").Value = "'0.800"
$ws.Range("E40").Value = "  -1.09%  "

$ws.Range("D41").Value = "'4.45"
$ws.Range("E41").Value = "  -2.99%  "

$ws.Range("D42").Value = "'6.22"
$ws.Range("E42").Value = "  -4.48%  "

$ws.Range("E43").Value = "  -1.19%  "

$ws.Range("E44").Value = "  -2.65%  "

$ws.Range("D45").Value = "'24.42"
$ws.Range("E45").Value = "  -3.72%  "

$ws.Range("D46").Value = "'2.35"
$ws.Range("E46").Value = "  -6.38%  "

$ws.Range("D47").Value = "'326.79"
$ws.Range("E47").Value = "  -2.46%  "

$ws.Range("E48").Value = "  -3.74%  "

$ws.Range("E49").Value = "  +1.15%  "

$ws.Range("D50").Value = "'0.979"
$ws.Range("E50").Value = "  -3.29%  "

$ws.Range("D51").Value = "'6.16"
$ws.Range("E51").Value = "  -2.28%  "
